$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Update the costs in column C (new dates and costs)
$ws.Range("C2").Value = 246
$ws.Range("C3").Value = 278.94
$ws.Range("C4").Value = 385.60000000000002

# Move the Global sheet's cursor to D7 without disturbing which sheet/tab
# is actually active in the workbook.
$origActive = $wb.ActiveSheet
$ws.Range("D7").Select()
$origActive.Activate()
